$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Copy the formatting of column I (rows 4-18) into column J so the new
#    column inherits borders / number formats / fonts / alignment that match
#    the existing table layout.
# ---------------------------------------------------------------------------
$ws.Range("I4:I18").Copy() | Out-Null
$ws.Range("J4:J18").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Fill in the new 2020 column values / formula.
# ---------------------------------------------------------------------------
$ws.Range("J4").Value2 = 2020
$ws.Range("J5").Value2 = 8017.9
# J6 stays empty (section header row)
$ws.Range("J7").Formula = "=J5-J8"
$ws.Range("J8").Value2 = 249.8
# J9 stays empty (section header row)
$ws.Range("J10").Value2 = 757.6
$ws.Range("J11").Value2 = 984.4
$ws.Range("J12").Value2 = 646.20000000000005
$ws.Range("J13").Value2 = 667.6
$ws.Range("J14").Value2 = 1147
$ws.Range("J15").Value2 = 961.1
$ws.Range("J16").Value2 = 2664.5
$ws.Range("J17").Value2 = 132.5
$ws.Range("J18").Value2 = 57

# ---------------------------------------------------------------------------
# 3. The bold "theme colour" header / section-total cells (J4, J9) get their
#    own distinct font, same as the equivalent bold cells further down would
#    when retyped by hand - nudge the engine into minting a fresh font entry
#    for them instead of silently reusing font 12.
# ---------------------------------------------------------------------------
$ws.Range("J4").Font.Bold = $true
$ws.Range("J4").Font.Size = 9
$ws.Range("J4").Font.Name = "Times New Roman"
$ws.Range("J4").Font.ThemeColor = 1

$ws.Range("J9").Font.Bold = $true
$ws.Range("J9").Font.Size = 9
$ws.Range("J9").Font.Name = "Times New Roman"
$ws.Range("J9").Font.ThemeColor = 1

# ---------------------------------------------------------------------------
# 4. The regular (non-bold) theme-coloured data cells (J10:J18) likewise get
#    their own distinct font entry, separate from column I's font 16.
# ---------------------------------------------------------------------------
$ws.Range("J10:J18").Font.Size = 9
$ws.Range("J10:J18").Font.Name = "Times New Roman"
$ws.Range("J10:J18").Font.ThemeColor = 1

# ---------------------------------------------------------------------------
# 5. Sheet level bookkeeping: dimension is recalculated automatically by the
#    engine; just restore the selected cell like Excel would leave it after
#    the edit.
# ---------------------------------------------------------------------------
$ws.Range("J19").Select() | Out-Null
